$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 641.5714
$ws.Range("I29").Value = 672.75
$ws.Range("K29").Value = 2018.25
$ws.Range("M29").Value = -1737.25

# Row 135
$ws.Range("H135").Value = 661.1667
$ws.Range("I135").Value = 698.5454999999999
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 6286.9095
$ws.Range("L135").Value = 2250
$ws.Range("M135").Value = -3751.9095
$ws.Range("N135").Value = -7320

# Row 137
$ws.Range("H137").Value = 3892.6924
$ws.Range("I137").Value = 4142.909
$ws.Range("J137").Value = 2516.5
$ws.Range("K137").Value = 12428.727
$ws.Range("L137").Value = 7549.5
$ws.Range("M137").Value = -9878.726999999999
$ws.Range("N137").Value = -12649.5

# Row 138
$ws.Range("H138").Value = 4250.543
$ws.Range("J138").Value = 5017.4443
$ws.Range("L138").Value = 15052.3329
$ws.Range("N138").Value = -25332.3329

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3270709.2
$ws.Range("I2").Value = 4557.2
$ws.Range("J2").Value = 7353399
$ws.Range("K2").Value = 4557.2
$ws.Range("L2").Value = 7353399
$ws.Range("M2").Value = -4444.2
$ws.Range("N2").Value = -7353625

# Row 5
$ws.Range("H5").Value = 195.5
$ws.Range("I5").Value = 90
$ws.Range("K5").Value = 90
$ws.Range("M5").Value = 22

# Row 45
$ws.Range("H45").Value = 85392.75
$ws.Range("I45").Value = 112866.555
$ws.Range("J45").Value = 2971.3333
$ws.Range("K45").Value = 112866.555
$ws.Range("L45").Value = 2971.3333
$ws.Range("M45").Value = -112489.555
$ws.Range("N45").Value = -3725.3333

# Row 74
$ws.Range("H74").Value = 5006191
$ws.Range("I74").Value = 10000621
$ws.Range("J74").Value = 11760.8
$ws.Range("K74").Value = 10000621
$ws.Range("L74").Value = 11760.8
$ws.Range("M74").Value = -9999747
$ws.Range("N74").Value = -13508.8

# Row 77
$ws.Range("H77").Value = 5006191
$ws.Range("I77").Value = 10000621
$ws.Range("J77").Value = 11760.8
$ws.Range("K77").Value = 50003105
$ws.Range("L77").Value = 58804
$ws.Range("M77").Value = -49998737
$ws.Range("N77").Value = -67540

# Row 116
$ws.Range("H116").Value = 3270709.2
$ws.Range("I116").Value = 4557.2
$ws.Range("J116").Value = 7353399
$ws.Range("K116").Value = 4557.2
$ws.Range("L116").Value = 7353399
$ws.Range("M116").Value = -2263.2
$ws.Range("N116").Value = -7357987

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3270709.2
$ws.Range("I3").Value = 4557.2
$ws.Range("J3").Value = 7353399
$ws.Range("K3").Value = 4557.2
$ws.Range("L3").Value = 7353399
$ws.Range("M3").Value = -4443.2
$ws.Range("N3").Value = -7353627

# Row 4
$ws.Range("H4").Value = 195.5
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 90
$ws.Range("M4").Value = 25

# Row 64
$ws.Range("H64").Value = 2728.9092
$ws.Range("I64").Value = 5651
$ws.Range("J64").Value = 1059.1428
$ws.Range("K64").Value = 5651
$ws.Range("L64").Value = 1059.1428
$ws.Range("M64").Value = -5426
$ws.Range("N64").Value = -1509.1428

# Row 67
$ws.Range("H67").Value = 2728.9092
$ws.Range("I67").Value = 5651
$ws.Range("J67").Value = 1059.1428
$ws.Range("K67").Value = 5651
$ws.Range("L67").Value = 1059.1428
$ws.Range("M67").Value = -4871
$ws.Range("N67").Value = -2619.1428

# Row 107
$ws.Range("H107").Value = 1799.1111
$ws.Range("I107").Value = 1491.421
$ws.Range("J107").Value = 2529.875
$ws.Range("K107").Value = 1491.421
$ws.Range("L107").Value = 2529.875
$ws.Range("M107").Value = 428.579
$ws.Range("N107").Value = -6369.875

$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 2032.3572
$ws.Range("I35").Value = 586.63635
$ws.Range("K35").Value = 586.63635
$ws.Range("M35").Value = -292.63635

# Row 132
$ws.Range("H132").Value = 3326.8386
$ws.Range("I132").Value = 3090.9412
$ws.Range("J132").Value = 3613.2856
$ws.Range("K132").Value = 9272.8236
$ws.Range("L132").Value = 10839.8568
$ws.Range("M132").Value = -6742.8236
$ws.Range("N132").Value = -15899.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 88.454544
$ws.Range("I40").Value = 57.3
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 229.2
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -160.2
$ws.Range("N40").Value = -1738

# Row 131
$ws.Range("H131").Value = 1001522.4
$ws.Range("I131").Value = 1557.375
$ws.Range("J131").Value = 1088475.9
$ws.Range("K131").Value = 4672.125
$ws.Range("L131").Value = 3265427.7
$ws.Range("M131").Value = 367.875
$ws.Range("N131").Value = -3275507.7

# Row 133
$ws.Range("H133").Value = 6393.5293
$ws.Range("I133").Value = 3733.3333
$ws.Range("J133").Value = 6963.5713
$ws.Range("K133").Value = 11199.9999
$ws.Range("L133").Value = 20890.7139
$ws.Range("M133").Value = -6139.999899999999
$ws.Range("N133").Value = -31010.7139

$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

# Row 82
$ws.Range("H82").Value = 1683.6666
$ws.Range("I82").Value = 1025.5
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1025.5
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -664.5
$ws.Range("N82").Value = -3722

# Row 85
$ws.Range("H85").Value = 1683.6666
$ws.Range("I85").Value = 1025.5
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1025.5
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = 222.5
$ws.Range("N85").Value = -5496

# Row 136
$ws.Range("H136").Value = 2727.4546
$ws.Range("I136").Value = 1778
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 5334
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -2784
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 786.9091
$ws.Range("I136").Value = 738.0625
$ws.Range("J136").Value = 832.8823
$ws.Range("K136").Value = 2214.1875
$ws.Range("L136").Value = 2498.6469
$ws.Range("M136").Value = 335.8125
$ws.Range("N136").Value = -7598.6469
